# Applies the two changes from the commit:
#   1. Change the table style used by the budgeting table on slide 16
#      from {C7266158-205F-4FEB-8C4A-A1C8C010DADF} to
#      {D53FD201-F5F4-414C-8F30-290557C0FAFB}.
#   2. Swap the deck's colour theme ("Integral") for the stock
#      "Office Theme" colour palette.

$p = $ppt.ActivePresentation

# Helper: VBA-style RGB() -> OLE/BGR packed integer used by ColorFormat.RGB
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Table style change (slide 16, 3rd shape, the 2-column table) ---
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{D53FD201-F5F4-414C-8F30-290557C0FAFB}")

# --- 2) Theme colour swap: "Integral" palette -> "Office Theme" palette ---
# Order matches MsoThemeColorSchemeIndex / the <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),   # dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),   # dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),   # accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),   # accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),   # accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),   # accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),   # hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
